$d = $word.ActiveDocument

# "Working Initial Prototype" bullet list -------------------------------

# 1) "Program the ball controller ..." bullet is done -> strike it out.
#    (Single run already; just flip StrikeThrough on for the paragraph.)
$p1 = $d.Paragraphs.Item(4)
$p1.Range.Font.StrikeThrough = $true

# 2) "Create prefabs ..." bullet is done too -> collapse its three runs
#    ("...to be ", "play tested", ".") into one run by replacing the
#    whole sentence with itself, then strike it out.
$d.Content.Find.Execute(
    "Create prefabs for the various game objects that will act as obstacles and implement them in a basic level format to be play tested.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Create prefabs for the various game objects that will act as obstacles and implement them in a basic level format to be play tested.",
    2)
$p2 = $d.Paragraphs.Item(5)
$p2.Range.Font.StrikeThrough = $true

# "Iterating on Initial Prototype" bullet list ---------------------------

# 3) "Finalize the tuning ..." bullet: merge its three runs
#    ("... influence the ", "behavioural", " physics ...") into a single
#    run - no formatting change, just a text/run-structure cleanup.
$d.Content.Find.Execute(
    "Finalize the tuning of the variable constants which influence the behavioural physics of both the ball hitting mechanics as well as the environmental aspects.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Finalize the tuning of the variable constants which influence the behavioural physics of both the ball hitting mechanics as well as the environmental aspects.",
    2)

$d.Save()
